$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on every cell we touch so numeric-looking
# strings (e.g. "311.71", "1.007") are preserved as text, matching the
# source inlineStr cells instead of being auto-converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D2").Value = '26.913.79'
$ws.Range("E2").Value = '  +2.21%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D3").Value = '1.811.07'
$ws.Range("E3").Value = '  +2.88%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.62%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D5").Value = '311.71'
$ws.Range("E5").Value = '  +2.72%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  +0.65%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4291'
$ws.Range("E7").Value = '  +0.19%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3691'
$ws.Range("E8").Value = '  +2.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07238'
$ws.Range("E9").Value = '  +2.72%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8644'
$ws.Range("E10").Value = '  +4.15%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D11").Value = '2.052.87'
$ws.Range("E11").Value = '  +17.41%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D12").Value = '21.21'
$ws.Range("E12").Value = '  +5.52%  '

$ws.Range("B13").NumberFormat = "@"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '5.397'
$ws.Range("E13").Value = '  +3.57%  '

$ws.Range("B14").NumberFormat = "@"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '6.619'
$ws.Range("E14").Value = '  +4.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06935'
$ws.Range("E15").Value = '  +2.24%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D16").Value = '80.83'
$ws.Range("E16").Value = '  +2.37%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D17").Value = '1.012'
$ws.Range("E17").Value = '  +0.85%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008921'
$ws.Range("E18").Value = '  +3.53%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D19").Value = '1.007'
$ws.Range("E19").Value = '  +0.58%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D21").Value = '26.951.38'
$ws.Range("E21").Value = '  +1.91%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D22").Value = '5.193'
$ws.Range("E22").Value = '  +4.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D23").Value = '10.94'
$ws.Range("E23").Value = '  -1.12%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D24").Value = '2.288.23'
$ws.Range("E24").Value = '  +16.25%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D25").Value = '154.07'
$ws.Range("E25").Value = '  +1.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D26").Value = '1.885'
$ws.Range("E26").Value = '  -1.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D27").Value = '18.32'
$ws.Range("E27").Value = '  +1.26%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D28").Value = '5.235'
$ws.Range("E28").Value = '  +4.17%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D29").Value = '1.921'
$ws.Range("E29").Value = '  +14.90%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D30").Value = '114.63'
$ws.Range("E30").Value = '  +0.43%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08953'
$ws.Range("E31").Value = '  +0.89%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7421'
$ws.Range("E32").Value = '  +2.81%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.80%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D34").Value = '4.430'
$ws.Range("E34").Value = '  +3.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D35").Value = '2.807'
$ws.Range("E35").Value = '  +3.63%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.81%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D37").Value = '1.123'
$ws.Range("E37").Value = '  +5.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05232'
$ws.Range("E38").Value = '  +2.71%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01922'
$ws.Range("E39").Value = '  +2.47%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5086'
$ws.Range("E40").Value = '  +4.05%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D41").Value = '2.748'
$ws.Range("E41").Value = '  +11.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1652'
$ws.Range("E42").Value = '  +3.30%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D43").Value = '6.489'
$ws.Range("E43").Value = '  +5.76%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D44").Value = '8.289'
$ws.Range("E44").Value = '  +3.57%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D45").Value = '107.45'
$ws.Range("E45").Value = '  +2.91%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D46").Value = '10.42'
$ws.Range("E46").Value = '  +4.41%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.77%  '

$ws.Range("B48").NumberFormat = "@"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '1.646'
$ws.Range("E48").Value = '  +5.13%  '

$ws.Range("B49").NumberFormat = "@"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.06274'
$ws.Range("E49").Value = '  +1.64%  '

$ws.Range("B50").NumberFormat = "@"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").Value = '0.4548'
$ws.Range("E50").Value = '  +1.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("D51").Value = '1.815'
$ws.Range("E51").Value = '  +6.18%  '
